# Actualización desde MV -datos-
# Updates the last existing data row (75) with revised figures and appends
# a new data row (76) for the quarter 01-04-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 75 with corrected values ---
$ws.Range("B75").Value = 336711
$ws.Range("F75").Value = 69044
$ws.Range("H75").Value = 68667
$ws.Range("I75").Value = -686
$ws.Range("J75").Value = 75228
$ws.Range("K75").Value = 4405
$ws.Range("L75").Value = 70823
$ws.Range("M75").Value = 184013
$ws.Range("N75").Value = 180866
$ws.Range("Q75").Value = 7129
$ws.Range("R75").Value = 19579
$ws.Range("S75").Value = 317132
$ws.Range("U75").Value = 19782
$ws.Range("V75").Value = 10982
$ws.Range("W75").Value = 8801
$ws.Range("X75").Value = 55560
$ws.Range("Z75").Value = 52105
$ws.Range("AA75").Value = 13461
$ws.Range("AB75").Value = 929
$ws.Range("AC75").Value = 12533
$ws.Range("AD75").Value = 216752
$ws.Range("AE75").Value = 151058
$ws.Range("AH75").Value = 11080

# --- Append new row 76 for quarter 01-04-2021 ---
# Force the date-like label to stay plain text (not get auto-converted to a
# date serial) while keeping the cell's style at the workbook default.
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").ClearFormats()
$ws.Range("B76").Value = 340074
$ws.Range("C76").Value = 1979
$ws.Range("D76").Value = 1974
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 73298
$ws.Range("G76").Value = 103
$ws.Range("H76").Value = 73195
$ws.Range("I76").Value = -552
$ws.Range("J76").Value = 74093
$ws.Range("K76").Value = 4733
$ws.Range("L76").Value = 69360
$ws.Range("M76").Value = 183804
$ws.Range("N76").Value = 180720
$ws.Range("O76").Value = 1172
$ws.Range("P76").Value = 1912
$ws.Range("Q76").Value = 7452
$ws.Range("R76").Value = 12280
$ws.Range("S76").Value = 327795
$ws.Range("T76").Value = 439
$ws.Range("U76").Value = 23867
$ws.Range("V76").Value = 12990
$ws.Range("W76").Value = 10877
$ws.Range("X76").Value = 56444
$ws.Range("Y76").Value = 508
$ws.Range("Z76").Value = 55935
$ws.Range("AA76").Value = 13643
$ws.Range("AB76").Value = 934
$ws.Range("AC76").Value = 12709
$ws.Range("AD76").Value = 221240
$ws.Range("AE76").Value = 156058
$ws.Range("AF76").Value = 2030
$ws.Range("AG76").Value = 63153
$ws.Range("AH76").Value = 12161
